# Append profile rows 6-10 (Careops 1: Basic UI and patient dashboard)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 6; $i -le 10; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
    $ws.Cells.Item($i, 2).Value = "UnifiedTestProfile$i"
    $ws.Cells.Item($i, 3).Value = "Available"
}
